# Update the "new installation" IBAMR_SRC_DIR / IBAMR_BUILD_DIR lines on the
# "Updating examples from UNC" slide with the new cluster path.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(23)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

# --- Paragraph 7: "IBAMR_SRC_DIR = " -------------------------------------
$para = $tr.Paragraphs(7, 1)
$oldText = $para.Text.TrimEnd([char]13)          # drop the trailing paragraph mark
$sub = $tr.Characters($para.Start, $oldText.Length)
$sub.Text = $oldText + "/groups/lauram9/ib10/"   # extend existing run in place
$para.InsertAfter("ibamr") | Out-Null
$para.InsertAfter("/IBAMR ") | Out-Null

# --- Paragraph 8: "IBAMR_BUILD_DIR = " ------------------------------------
# re-fetch the TextRange/paragraph since prior inserts shifted character offsets
$tr = $tf.TextRange
$para = $tr.Paragraphs(8, 1)
$oldText = $para.Text.TrimEnd([char]13)
$sub = $tr.Characters($para.Start, $oldText.Length)
$sub.Text = $oldText + "/groups/lauram9/ib10/"
$para.InsertAfter("ibamr") | Out-Null
$para.InsertAfter("/") | Out-Null
$para.InsertAfter("ibamr") | Out-Null
$para.InsertAfter("-") | Out-Null
$para.InsertAfter("objs") | Out-Null
$para.InsertAfter("-opt") | Out-Null
